$d = $word.ActiveDocument

# Update the date heading (first paragraph) in place, preserving its run formatting.
$dateRange = $d.Paragraphs.Item(1).Range
$dateRange.Text = "2025-11-22 Saturday"

# Update each arithmetic-problem cell directly via Table/Row/Column addressing and
# Range.Text assignment, so a replacement can never bleed into a neighboring cell
# (several of the expressions are substrings of one another, e.g. "6+48=" inside "46+48=").
$tbl = $d.Tables.Item(1)

function Set-CellText($row, $col, $new) {
    $tbl.Cell($row, $col).Range.Text = $new
}

Set-CellText 1 1 "92-17="  # was "55+8="
Set-CellText 1 2 "63+28="  # was "77-18="
Set-CellText 1 3 "26+39="  # was "6+48="
Set-CellText 1 4 "50-44="  # was "25+49="
Set-CellText 1 5 "8+88="  # was "73-55="
Set-CellText 2 1 "40-32="  # was "39+25="
Set-CellText 2 2 "27+59="  # was "53-46="
Set-CellText 2 3 "33+38="  # was "55+37="
Set-CellText 2 4 "39+9="  # was "16+15="
Set-CellText 2 5 "59+35="  # was "84-15="
Set-CellText 3 1 "8+65="  # was "60-32="
Set-CellText 3 2 "81-24="  # was "54-27="
Set-CellText 3 3 "20-3="  # was "38+49="
Set-CellText 3 4 "7+87="  # was "49+12="
Set-CellText 3 5 "90-11="  # was "27+45="
Set-CellText 4 1 "78+3="  # was "53-4="
Set-CellText 4 2 "75+7="  # was "90-33="
Set-CellText 4 3 "69+19="  # was "75-6="
Set-CellText 4 4 "34-6="  # was "29+43="
Set-CellText 4 5 "71-36="  # was "46-38="
Set-CellText 5 1 "60-1="  # was "81-37="
Set-CellText 5 2 "29+8="  # was "64+7="
Set-CellText 5 3 "93-85="  # was "55+9="
Set-CellText 5 4 "76+19="  # was "56+7="
Set-CellText 5 5 "77-68="  # was "62-55="
Set-CellText 6 1 "55-9="  # was "82-45="
Set-CellText 6 2 "22-9="  # was "26+16="
Set-CellText 6 3 "40-22="  # was "37+34="
Set-CellText 6 4 "19+53="  # was "24+18="
Set-CellText 6 5 "8+67="  # was "14+59="
Set-CellText 7 1 "19+15="  # was "4+57="
Set-CellText 7 2 "93-18="  # was "63-25="
Set-CellText 7 3 "43-5="  # was "44-38="
Set-CellText 7 4 "29+52="  # was "76-59="
Set-CellText 7 5 "31-5="  # was "20-12="
Set-CellText 8 1 "18+33="  # was "48+13="
Set-CellText 8 2 "92-63="  # was "9+88="
Set-CellText 8 3 "37-18="  # was "86-59="
Set-CellText 8 4 "62-18="  # was "17+54="
Set-CellText 8 5 "19+68="  # was "72-63="
Set-CellText 9 1 "5+18="  # was "67+5="
Set-CellText 9 2 "59+6="  # was "22+39="
Set-CellText 9 3 "63-24="  # was "28+29="
Set-CellText 9 4 "19+3="  # was "77+17="
Set-CellText 9 5 "8+33="  # was "9+58="
Set-CellText 10 1 "28+56="  # was "84-67="
Set-CellText 10 2 "4+58="  # was "21-17="
Set-CellText 10 3 "38+35="  # was "29+33="
Set-CellText 10 4 "70-5="  # was "60-49="
Set-CellText 10 5 "31-9="  # was "25+18="
Set-CellText 11 1 "42-5="  # was "36+55="
Set-CellText 11 2 "70-56="  # was "71-27="
Set-CellText 11 3 "2+59="  # was "9+89="
Set-CellText 11 4 "42+39="  # was "8+79="
Set-CellText 11 5 "17+49="  # was "36+35="
Set-CellText 12 1 "74+8="  # was "55-38="
Set-CellText 12 2 "14+67="  # was "92-64="
Set-CellText 12 3 "43-26="  # was "90-46="
Set-CellText 12 4 "60-9="  # was "74-35="
Set-CellText 12 5 "73-66="  # was "3+68="
Set-CellText 13 1 "96-7="  # was "84-27="
Set-CellText 13 2 "9+25="  # was "9+87="
Set-CellText 13 3 "11-2="  # was "39+12="
Set-CellText 13 4 "27+56="  # was "61-33="
Set-CellText 13 5 "92-27="  # was "9+39="
Set-CellText 14 1 "58+15="  # was "43-8="
Set-CellText 14 2 "93-16="  # was "90-53="
Set-CellText 14 3 "52-49="  # was "14+27="
Set-CellText 14 4 "80-78="  # was "46+48="
Set-CellText 14 5 "86+7="  # was "4+77="
Set-CellText 15 1 "72-35="  # was "39+24="
Set-CellText 15 2 "66-47="  # was "52-13="
Set-CellText 15 3 "40-7="  # was "38+4="
Set-CellText 15 4 "91-3="  # was "19+19="
Set-CellText 15 5 "41-7="  # was "34+8="
Set-CellText 16 1 "15+78="  # was "23-14="
Set-CellText 16 2 "29+19="  # was "44+27="
Set-CellText 16 3 "86-39="  # was "9+57="
Set-CellText 16 4 "91-16="  # was "5+36="
Set-CellText 16 5 "71-59="  # was "45+19="
Set-CellText 17 1 "11-2="  # was "6+47="
Set-CellText 17 2 "59+36="  # was "89+2="
Set-CellText 17 3 "78+3="  # was "40-15="
Set-CellText 17 4 "37+55="  # was "25+9="
Set-CellText 17 5 "62-58="  # was "48+28="
Set-CellText 18 1 "54-48="  # was "26+28="
Set-CellText 18 2 "7+34="  # was "28+58="
Set-CellText 18 3 "51-12="  # was "4+47="
Set-CellText 18 4 "18+3="  # was "85-69="
Set-CellText 18 5 "61-36="  # was "72-58="
Set-CellText 19 1 "60-27="  # was "29+22="
Set-CellText 19 2 "74-17="  # was "63-26="
Set-CellText 19 3 "82-25="  # was "9+14="
Set-CellText 19 4 "72-14="  # was "28+49="
Set-CellText 19 5 "92-4="  # was "9+56="
Set-CellText 20 1 "47+18="  # was "8+66="
Set-CellText 20 2 "13+58="  # was "56+8="
Set-CellText 20 3 "71-5="  # was "70-6="
Set-CellText 20 4 "57-49="  # was "2+39="
Set-CellText 20 5 "70-31="  # was "20-19="
